$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 57-58, shifting existing rows 57:142 down to 59:144
# (this also grows the sheet's used-range dimension from R142 to R144).
$ws.Rows("57:58").Insert()

# Fill the newly inserted row 57 with its data.
$ws.Range("A57").Value = 6
$ws.Range("B57").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C57").Value = "Metropolitana"
$ws.Range("D57").Value = 44482
$ws.Range("E57").Value = 13
$ws.Range("F57").Value = 100112022
$ws.Range("G57").Value = "Arveja Verde"
$ws.Range("H57").Value = "Perfection"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 120
$ws.Range("K57").Value = 22000
$ws.Range("L57").Value = 23000
$ws.Range("M57").Value = 22417
$ws.Range("N57").Value = "$/malla 25 kilos"
$ws.Range("O57").Value = "Provincia de Huasco"
$ws.Range("P57").Value = 897
$ws.Range("Q57").Value = 25
$ws.Range("R57").Value = "Hortaliza"

# Fill the newly inserted row 58 with its data.
$ws.Range("A58").Value = 6
$ws.Range("B58").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C58").Value = "Metropolitana"
$ws.Range("D58").Value = 44482
$ws.Range("E58").Value = 13
$ws.Range("F58").Value = 100112022
$ws.Range("G58").Value = "Arveja Verde"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 220
$ws.Range("K58").Value = 18000
$ws.Range("L58").Value = 20000
$ws.Range("M58").Value = 18909
$ws.Range("N58").Value = "$/saco 25 kilos"
$ws.Range("O58").Value = "Región Metropolitana"
$ws.Range("P58").Value = 756
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = "Hortaliza"
